$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 316, shifting existing rows 316:357 down to 317:358.
$ws.Rows("316:316").Insert()

# Populate the newly inserted row 316 with its data.
$ws.Range("A316").Value = 6
$ws.Range("B316").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C316").Value = "Metropolitana"
$ws.Range("D316").Value = 45218
$ws.Range("E316").Value = 13
$ws.Range("F316").Value = 100112022
$ws.Range("G316").Value = "Arveja Verde"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 300
$ws.Range("K316").Value = 18000
$ws.Range("L316").Value = 20000
$ws.Range("M316").Value = 18800
$ws.Range("N316").Value = "$/saco 25 kilos"
$ws.Range("O316").Value = "Región de O'Higgins"
$ws.Range("P316").Value = 752
$ws.Range("Q316").Value = 25
$ws.Range("R316").Value = "Hortaliza"
